# Auto-generated updates for Sheets/Lich_Profits.xlsx (per-sheet tabs ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR)
# Commit: chore: update Sheets via scheduled runner
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(12, 8).Value = 127.2  # H12: was 107.5
$ws.Cells.Item(12, 9).Value = 155  # I12: was 125.8
$ws.Cells.Item(12, 11).Value = 155  # K12: was 125.8
$ws.Cells.Item(12, 13).Value = 15  # M12: was 44.2
$ws.Cells.Item(32, 8).Value = 2999  # H32: was 2999.5
$ws.Cells.Item(32, 9).Value = 2999  # I32: was 2999.5
$ws.Cells.Item(32, 11).Value = 2999  # K32: was 2999.5
$ws.Cells.Item(32, 13).Value = -2673  # M32: was -2673.5
$ws.Cells.Item(51, 8).Value = 9820.893  # H51: was 9892.296
$ws.Cells.Item(51, 9).Value = 16331.167  # I51: was 16331.333
$ws.Cells.Item(51, 10).Value = 8045.364  # J51: was 8052.5713
$ws.Cells.Item(51, 11).Value = 16331.167  # K51: was 16331.333
$ws.Cells.Item(51, 12).Value = 8045.364  # L51: was 8052.5713
$ws.Cells.Item(51, 13).Value = -15847.167  # M51: was -15847.333
$ws.Cells.Item(51, 14).Value = -9013.364  # N51: was -9020.5713
$ws.Cells.Item(88, 8).Value = 2955  # H88: was 3174.3635
$ws.Cells.Item(88, 9).Value = 2739  # I88: was 3048
$ws.Cells.Item(88, 10).Value = 3090  # J88: was 3246.5715
$ws.Cells.Item(88, 11).Value = 2739  # K88: was 3048
$ws.Cells.Item(88, 12).Value = 3090  # L88: was 3246.5715
$ws.Cells.Item(88, 13).Value = -2333  # M88: was -2642
$ws.Cells.Item(88, 14).Value = -3902  # N88: was -4058.5715
$ws.Cells.Item(91, 8).Value = 2955  # H91: was 3174.3635
$ws.Cells.Item(91, 9).Value = 2739  # I91: was 3048
$ws.Cells.Item(91, 10).Value = 3090  # J91: was 3246.5715
$ws.Cells.Item(91, 11).Value = 2739  # K91: was 3048
$ws.Cells.Item(91, 12).Value = 3090  # L91: was 3246.5715
$ws.Cells.Item(91, 13).Value = -1335  # M91: was -1644
$ws.Cells.Item(91, 14).Value = -5898  # N91: was -6054.5715
$ws.Cells.Item(103, 8).Value = 758.03845  # H103: was 683.8125
$ws.Cells.Item(103, 9).Value = 570.6429000000001  # I103: was 508.1
$ws.Cells.Item(103, 11).Value = 1711.9287  # K103: was 1524.3
$ws.Cells.Item(103, 13).Value = -1125.9287  # M103: was -938.3000000000002
$ws.Cells.Item(121, 8).Value = 5665.6665  # H121: was 5874.25
$ws.Cells.Item(121, 10).Value = 5665.6665  # J121: was 5874.25
$ws.Cells.Item(121, 12).Value = 16996.9995  # L121: was 17622.75
$ws.Cells.Item(121, 14).Value = -20490.9995  # N121: was -21116.75
$ws.Cells.Item(137, 8).Value = 39007.87  # H137: was 37790.75
$ws.Cells.Item(137, 10).Value = 7049.2  # J137: was 6413.8184
$ws.Cells.Item(137, 12).Value = 21147.6  # L137: was 19241.4552
$ws.Cells.Item(137, 14).Value = -26247.6  # N137: was -24341.4552
$ws.Cells.Item(138, 8).Value = 3377.1707  # H138: was 3422.2563
$ws.Cells.Item(138, 10).Value = 3857.5  # J138: was 3948.1333
$ws.Cells.Item(138, 12).Value = 11572.5  # L138: was 11844.3999
$ws.Cells.Item(138, 14).Value = -21852.5  # N138: was -22124.3999

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 2845.9092  # H2: was 3061
$ws.Cells.Item(2, 9).Value = 2310.4443  # I2: was 2512.375
$ws.Cells.Item(2, 11).Value = 2310.4443  # K2: was 2512.375
$ws.Cells.Item(2, 13).Value = -2197.4443  # M2: was -2399.375
$ws.Cells.Item(48, 8).Value = 295000  # H48: was 0
$ws.Cells.Item(48, 10).Value = 295000  # J48: was 0
$ws.Cells.Item(48, 12).Value = 295000  # L48: was 0
$ws.Cells.Item(48, 14).Value = -295768  # N48: was None
$ws.Cells.Item(63, 8).Value = 3893.6  # H63: was 4022.3572
$ws.Cells.Item(63, 9).Value = 3472.2727  # I63: was 3622
$ws.Cells.Item(63, 10).Value = 5052.25  # J63: was 4743
$ws.Cells.Item(63, 11).Value = 3472.2727  # K63: was 3622
$ws.Cells.Item(63, 12).Value = 5052.25  # L63: was 4743
$ws.Cells.Item(63, 13).Value = -2786.2727  # M63: was -2936
$ws.Cells.Item(63, 14).Value = -6424.25  # N63: was -6115
$ws.Cells.Item(66, 8).Value = 3893.6  # H66: was 4022.3572
$ws.Cells.Item(66, 9).Value = 3472.2727  # I66: was 3622
$ws.Cells.Item(66, 10).Value = 5052.25  # J66: was 4743
$ws.Cells.Item(66, 11).Value = 17361.3635  # K66: was 18110
$ws.Cells.Item(66, 12).Value = 25261.25  # L66: was 23715
$ws.Cells.Item(66, 13).Value = -13929.3635  # M66: was -14678
$ws.Cells.Item(66, 14).Value = -32125.25  # N66: was -30579
$ws.Cells.Item(88, 8).Value = 1945.1666  # H88: was 2013.5834
$ws.Cells.Item(88, 9).Value = 1244  # I88: was 1380.8334
$ws.Cells.Item(88, 11).Value = 1244  # K88: was 1380.8334
$ws.Cells.Item(88, 13).Value = -838  # M88: was -974.8334
$ws.Cells.Item(91, 8).Value = 1945.1666  # H91: was 2013.5834
$ws.Cells.Item(91, 9).Value = 1244  # I91: was 1380.8334
$ws.Cells.Item(91, 11).Value = 1244  # K91: was 1380.8334
$ws.Cells.Item(91, 13).Value = 160  # M91: was 23.16660000000002
$ws.Cells.Item(97, 8).Value = 1454.1538  # H97: was 1600.2858
$ws.Cells.Item(97, 10).Value = 3749.5  # J97: was 3666.3333
$ws.Cells.Item(97, 12).Value = 3749.5  # L97: was 3666.3333
$ws.Cells.Item(97, 14).Value = -4741.5  # N97: was -4658.3333
$ws.Cells.Item(116, 8).Value = 2845.9092  # H116: was 3061
$ws.Cells.Item(116, 9).Value = 2310.4443  # I116: was 2512.375
$ws.Cells.Item(116, 11).Value = 2310.4443  # K116: was 2512.375
$ws.Cells.Item(116, 13).Value = -16.44430000000011  # M116: was -218.375
$ws.Cells.Item(122, 8).Value = 2337.28  # H122: was 2252.75
$ws.Cells.Item(122, 9).Value = 2347.375  # I122: was 2344.52
$ws.Cells.Item(122, 10).Value = 2095  # J122: was 1488
$ws.Cells.Item(122, 11).Value = 7042.125  # K122: was 7033.559999999999
$ws.Cells.Item(122, 12).Value = 6285  # L122: was 4464
$ws.Cells.Item(122, 13).Value = -4592.125  # M122: was -4583.559999999999
$ws.Cells.Item(122, 14).Value = -11185  # N122: was -9364

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 2845.9092  # H3: was 3061
$ws.Cells.Item(3, 9).Value = 2310.4443  # I3: was 2512.375
$ws.Cells.Item(3, 11).Value = 2310.4443  # K3: was 2512.375
$ws.Cells.Item(3, 13).Value = -2196.4443  # M3: was -2398.375
$ws.Cells.Item(22, 8).Value = 364  # H22: was 386.125
$ws.Cells.Item(22, 9).Value = 364  # I22: was 386.125
$ws.Cells.Item(22, 11).Value = 364  # K22: was 386.125
$ws.Cells.Item(22, 13).Value = -191  # M22: was -213.125
$ws.Cells.Item(82, 8).Value = 111154250  # H82: was 125047410
$ws.Cells.Item(82, 9).Value = 200013650  # I82: was 250014820
$ws.Cells.Item(82, 11).Value = 200013650  # K82: was 250014820
$ws.Cells.Item(82, 13).Value = -200013267  # M82: was -250014437
$ws.Cells.Item(85, 8).Value = 111154250  # H85: was 125047410
$ws.Cells.Item(85, 9).Value = 200013650  # I85: was 250014820
$ws.Cells.Item(85, 11).Value = 200013650  # K85: was 250014820
$ws.Cells.Item(85, 13).Value = -200012324  # M85: was -250013494
$ws.Cells.Item(107, 8).Value = 3620.9524  # H107: was 3659.5
$ws.Cells.Item(107, 9).Value = 3715.125  # I107: was 3772.8
$ws.Cells.Item(107, 11).Value = 3715.125  # K107: was 3772.8
$ws.Cells.Item(107, 13).Value = -1795.125  # M107: was -1852.8
$ws.Cells.Item(134, 8).Value = 2688.8667  # H134: was 2705.4138
$ws.Cells.Item(134, 9).Value = 2383.5715  # I134: was 2412.6667
$ws.Cells.Item(134, 11).Value = 7150.7145  # K134: was 7238.000100000001
$ws.Cells.Item(134, 13).Value = -4615.7145  # M134: was -4703.000100000001

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(22, 8).Value = 1112  # H22: was 1112.3334
$ws.Cells.Item(22, 9).Value = 420  # I22: was 418.5
$ws.Cells.Item(22, 10).Value = 1458  # J22: was 2500
$ws.Cells.Item(22, 11).Value = 420  # K22: was 418.5
$ws.Cells.Item(22, 12).Value = 1458  # L22: was 2500
$ws.Cells.Item(22, 13).Value = -70  # M22: was -68.5
$ws.Cells.Item(22, 14).Value = -2158  # N22: was -3200
$ws.Cells.Item(31, 8).Value = 234572.72  # H31: was 245981.14
$ws.Cells.Item(31, 9).Value = 477434.94  # I31: was 501271.7
$ws.Cells.Item(31, 10).Value = 2749.682  # J31: was 2847.2856
$ws.Cells.Item(31, 11).Value = 477434.94  # K31: was 501271.7
$ws.Cells.Item(31, 12).Value = 2749.682  # L31: was 2847.2856
$ws.Cells.Item(31, 13).Value = -477139.94  # M31: was -500976.7
$ws.Cells.Item(31, 14).Value = -3339.682  # N31: was -3437.2856
$ws.Cells.Item(34, 8).Value = 234572.72  # H34: was 245981.14
$ws.Cells.Item(34, 9).Value = 477434.94  # I34: was 501271.7
$ws.Cells.Item(34, 10).Value = 2749.682  # J34: was 2847.2856
$ws.Cells.Item(34, 11).Value = 477434.94  # K34: was 501271.7
$ws.Cells.Item(34, 12).Value = 2749.682  # L34: was 2847.2856
$ws.Cells.Item(34, 13).Value = -477232.94  # M34: was -501069.7
$ws.Cells.Item(34, 14).Value = -3153.682  # N34: was -3251.2856
$ws.Cells.Item(58, 8).Value = 4060.1853  # H58: was 4000.9
$ws.Cells.Item(58, 9).Value = 3819.0527  # I58: was 3758.2
$ws.Cells.Item(58, 10).Value = 4632.875  # J58: was 4486.3
$ws.Cells.Item(58, 11).Value = 3819.0527  # K58: was 3758.2
$ws.Cells.Item(58, 12).Value = 4632.875  # L58: was 4486.3
$ws.Cells.Item(58, 13).Value = -3616.0527  # M58: was -3555.2
$ws.Cells.Item(58, 14).Value = -5038.875  # N58: was -4892.3
$ws.Cells.Item(107, 8).Value = 4474.974  # H107: was 4701.108
$ws.Cells.Item(107, 9).Value = 689  # I107: was 761.2727
$ws.Cells.Item(107, 11).Value = 689  # K107: was 761.2727
$ws.Cells.Item(107, 13).Value = 1231  # M107: was 1158.7273
$ws.Cells.Item(132, 8).Value = 9079.857  # H132: was 10243.167
$ws.Cells.Item(132, 9).Value = 2190.75  # I132: was 2221
$ws.Cells.Item(132, 11).Value = 6572.25  # K132: was 6663
$ws.Cells.Item(132, 13).Value = -4042.25  # M132: was -4133
$ws.Cells.Item(136, 8).Value = 4060.1853  # H136: was 4000.9
$ws.Cells.Item(136, 9).Value = 3819.0527  # I136: was 3758.2
$ws.Cells.Item(136, 10).Value = 4632.875  # J136: was 4486.3
$ws.Cells.Item(136, 11).Value = 11457.1581  # K136: was 11274.6
$ws.Cells.Item(136, 12).Value = 13898.625  # L136: was 13458.9
$ws.Cells.Item(136, 13).Value = -8907.158100000001  # M136: was -8724.599999999999
$ws.Cells.Item(136, 14).Value = -18998.625  # N136: was -18558.9

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(21, 8).Value = 325  # H21: was 150
$ws.Cells.Item(21, 10).Value = 325  # J21: was 150
$ws.Cells.Item(21, 12).Value = 975  # L21: was 450
$ws.Cells.Item(21, 14).Value = -1321  # N21: was -796
$ws.Cells.Item(139, 8).Value = 2969.9473  # H139: was 2270.3684
$ws.Cells.Item(139, 9).Value = 2152.5833  # I139: was 1044.9166
$ws.Cells.Item(139, 11).Value = 6457.749899999999  # K139: was 3134.7498
$ws.Cells.Item(139, 13).Value = -1317.749899999999  # M139: was 2005.2502
$ws.Cells.Item(141, 8).Value = 4742.9  # H141: was 3992.5
$ws.Cells.Item(141, 9).Value = 3985.8  # I141: was 3992.5
$ws.Cells.Item(141, 10).Value = 5500  # J141: was 0
$ws.Cells.Item(141, 11).Value = 11957.4  # K141: was 11977.5
$ws.Cells.Item(141, 12).Value = 16500  # L141: was 0
$ws.Cells.Item(141, 13).Value = -6777.400000000001  # M141: was -6797.5
$ws.Cells.Item(141, 14).Value = -26860  # N141: was None

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(2, 8).Value = 350.5  # H2: was 270.1
$ws.Cells.Item(2, 9).Value = 382.6  # I2: was 295.125
$ws.Cells.Item(2, 10).Value = 190  # J2: was 170
$ws.Cells.Item(2, 11).Value = 382.6  # K2: was 295.125
$ws.Cells.Item(2, 12).Value = 190  # L2: was 170
$ws.Cells.Item(2, 13).Value = -269.6  # M2: was -182.125
$ws.Cells.Item(2, 14).Value = -416  # N2: was -396
$ws.Cells.Item(80, 8).Value = 4189.108  # H80: was 4219.1353
$ws.Cells.Item(80, 10).Value = 4546.9  # J80: was 4602.45
$ws.Cells.Item(80, 12).Value = 4546.9  # L80: was 4602.45
$ws.Cells.Item(80, 14).Value = -6542.9  # N80: was -6598.45
$ws.Cells.Item(83, 8).Value = 4189.108  # H83: was 4219.1353
$ws.Cells.Item(83, 10).Value = 4546.9  # J83: was 4602.45
$ws.Cells.Item(83, 12).Value = 22734.5  # L83: was 23012.25
$ws.Cells.Item(83, 14).Value = -32718.5  # N83: was -32996.25
$ws.Cells.Item(122, 8).Value = 6825.778  # H122: was 6589.8687
$ws.Cells.Item(122, 9).Value = 5293.1875  # I122: was 5081.1763
$ws.Cells.Item(122, 10).Value = 8051.85  # J122: was 7811.1904
$ws.Cells.Item(122, 11).Value = 15879.5625  # K122: was 15243.5289
$ws.Cells.Item(122, 12).Value = 24155.55  # L122: was 23433.5712
$ws.Cells.Item(122, 13).Value = -13429.5625  # M122: was -12793.5289
$ws.Cells.Item(122, 14).Value = -29055.55  # N122: was -28333.5712
$ws.Cells.Item(126, 8).Value = 9263  # H126: was 9870
$ws.Cells.Item(126, 9).Value = 10551.25  # I126: was 11515.714
$ws.Cells.Item(126, 11).Value = 31653.75  # K126: was 34547.142
$ws.Cells.Item(126, 13).Value = -29183.75  # M126: was -32077.142
$ws.Cells.Item(140, 8).Value = 70000  # H140: was 0
$ws.Cells.Item(140, 10).Value = 70000  # J140: was 0
$ws.Cells.Item(140, 12).Value = 70000  # L140: was 0
$ws.Cells.Item(140, 14).Value = -80360  # N140: was None

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(122, 8).Value = 3372.5  # H122: was 3458.2856
$ws.Cells.Item(122, 9).Value = 3347.3076  # I122: was 3418.4167
$ws.Cells.Item(122, 10).Value = 3700  # J122: was 3697.5
$ws.Cells.Item(122, 11).Value = 10041.9228  # K122: was 10255.2501
$ws.Cells.Item(122, 12).Value = 11100  # L122: was 11092.5
$ws.Cells.Item(122, 13).Value = -7591.9228  # M122: was -7805.250100000001
$ws.Cells.Item(122, 14).Value = -16000  # N122: was -15992.5

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(32, 8).Value = 0  # H32: was 2000
$ws.Cells.Item(32, 9).Value = 0  # I32: was 2000
$ws.Cells.Item(32, 11).Value = 0  # K32: was 2000
$ws.Cells.Item(32, 13).ClearContents()  # M32: was -1683, now empty
$ws.Cells.Item(51, 8).Value = 13041.083  # H51: was 11280.952
$ws.Cells.Item(51, 9).Value = 13041.083  # I51: was 13394.467
$ws.Cells.Item(51, 10).Value = 0  # J51: was 5997.1665
$ws.Cells.Item(51, 11).Value = 13041.083  # K51: was 13394.467
$ws.Cells.Item(51, 12).Value = 0  # L51: was 5997.1665
$ws.Cells.Item(51, 13).Value = -12531.083  # M51: was -12884.467
$ws.Cells.Item(51, 14).ClearContents()  # N51: was -7017.1665, now empty
$ws.Cells.Item(107, 8).Value = 435.7  # H107: was 435.8
$ws.Cells.Item(107, 9).Value = 435.7  # I107: was 435.8
$ws.Cells.Item(107, 11).Value = 1307.1  # K107: was 1307.4
$ws.Cells.Item(107, 13).Value = 612.9000000000001  # M107: was 612.5999999999999
$ws.Cells.Item(136, 8).Value = 1115489.1  # H136: was 772845.7
$ws.Cells.Item(136, 9).Value = 1669067.1  # I136: was 1002199.5
$ws.Cells.Item(136, 11).Value = 5007201.300000001  # K136: was 3006598.5
$ws.Cells.Item(136, 13).Value = -5004651.300000001  # M136: was -3004048.5

Write-Host "Applied all Lich_Profits updates"
